$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the item name for item 1 (Tea -> Green Tea)
$ws.Range("B2").Value = "Green Tea"

# Make sure Price for item 1 stays a number (10)
$ws.Range("C2").Value = 10

# Add a new item row: Item ID 10, Roti, Price 2, Special Item True, Is Active True
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Roti"
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = $true
$ws.Range("E11").Value = $true

# Widen column B to fit the new, longer item names
$ws.Columns.Item(2).ColumnWidth = 14.72
